$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pre-seed brand-new shared strings in the exact order they must be
# minted (2 chance, then 2*, then 1*, then 4*) so the saved sharedStrings
# table lands with indices 24-27 in that order, matching the target file.
$ws.Cells.Item(90, 1).Value = "2 chance"
$ws.Cells.Item(93, 5).Value = "2*"
$ws.Cells.Item(92, 6).Value = "1*"
$ws.Cells.Item(94, 8).Value = "4*"

# ---- Block: FIFO (header row 74) ----
$ws.Range("A59:N64").Copy()
$ws.Range("A75").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(74, 1).Value = "FIFO"
$ws.Cells.Item(75, 1).Value = "Marco/Pagina"
$ws.Cells.Item(75, 2).Value = 1
$ws.Cells.Item(75, 3).Value = 2
$ws.Cells.Item(75, 4).Value = 4
$ws.Cells.Item(75, 5).Value = 2
$ws.Cells.Item(75, 6).Value = 1
$ws.Cells.Item(75, 7).Value = 3
$ws.Cells.Item(75, 8).Value = 4
$ws.Cells.Item(75, 9).Value = 5
$ws.Cells.Item(75, 10).Value = 1
$ws.Cells.Item(75, 11).Value = 6
$ws.Cells.Item(75, 12).Value = 1
$ws.Cells.Item(75, 13).Value = 2
$ws.Cells.Item(75, 14).Value = 3
$ws.Cells.Item(76, 1).Value = "M1"
$ws.Cells.Item(76, 2).Value = 1
$ws.Cells.Item(76, 3).Value = 1
$ws.Cells.Item(76, 4).Value = 1
$ws.Cells.Item(76, 5).Value = 1
$ws.Cells.Item(76, 6).Value = 1
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 1
$ws.Cells.Item(76, 9).Value = 5
$ws.Cells.Item(77, 1).Value = "M2"
$ws.Cells.Item(77, 3).Value = 2
$ws.Cells.Item(77, 4).Value = 2
$ws.Cells.Item(77, 5).Value = 2
$ws.Cells.Item(77, 6).Value = 2
$ws.Cells.Item(77, 7).Value = 2
$ws.Cells.Item(77, 8).Value = 2
$ws.Cells.Item(77, 9).Value = 2
$ws.Cells.Item(78, 1).Value = "M3"
$ws.Cells.Item(78, 4).Value = 4
$ws.Cells.Item(78, 5).Value = 4
$ws.Cells.Item(78, 6).Value = 4
$ws.Cells.Item(78, 7).Value = 4
$ws.Cells.Item(78, 8).Value = 4
$ws.Cells.Item(78, 9).Value = 4
$ws.Cells.Item(79, 1).Value = "M4"
$ws.Cells.Item(79, 7).Value = 3
$ws.Cells.Item(79, 8).Value = 3
$ws.Cells.Item(79, 9).Value = 3
$ws.Cells.Item(80, 1).Value = "PF"
$ws.Cells.Item(80, 2).Value = "X"
$ws.Cells.Item(80, 3).Value = "X"
$ws.Cells.Item(80, 4).Value = "X"
$ws.Cells.Item(80, 7).Value = "X"
$ws.Cells.Item(80, 9).Value = "X"

# ---- Block: LRU (header row 82) ----
$ws.Range("A59:N64").Copy()
$ws.Range("A83").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(82, 1).Value = "LRU"
$ws.Cells.Item(83, 1).Value = "Marco/Pagina"
$ws.Cells.Item(83, 2).Value = 1
$ws.Cells.Item(83, 3).Value = 2
$ws.Cells.Item(83, 4).Value = 4
$ws.Cells.Item(83, 5).Value = 2
$ws.Cells.Item(83, 6).Value = 1
$ws.Cells.Item(83, 7).Value = 3
$ws.Cells.Item(83, 8).Value = 4
$ws.Cells.Item(83, 9).Value = 5
$ws.Cells.Item(83, 10).Value = 1
$ws.Cells.Item(83, 11).Value = 6
$ws.Cells.Item(83, 12).Value = 1
$ws.Cells.Item(83, 13).Value = 2
$ws.Cells.Item(83, 14).Value = 3
$ws.Cells.Item(84, 1).Value = "M1"
$ws.Cells.Item(84, 2).Value = 1
$ws.Cells.Item(84, 3).Value = 1
$ws.Cells.Item(84, 4).Value = 1
$ws.Cells.Item(84, 5).Value = 1
$ws.Cells.Item(84, 6).Value = 1
$ws.Cells.Item(84, 7).Value = 1
$ws.Cells.Item(84, 8).Value = 1
$ws.Cells.Item(84, 9).Value = 1
$ws.Cells.Item(85, 1).Value = "M2"
$ws.Cells.Item(85, 3).Value = 2
$ws.Cells.Item(85, 4).Value = 2
$ws.Cells.Item(85, 5).Value = 2
$ws.Cells.Item(85, 6).Value = 2
$ws.Cells.Item(85, 7).Value = 2
$ws.Cells.Item(85, 8).Value = 2
$ws.Cells.Item(85, 9).Value = 5
$ws.Cells.Item(86, 1).Value = "M3"
$ws.Cells.Item(86, 4).Value = 4
$ws.Cells.Item(86, 5).Value = 4
$ws.Cells.Item(86, 6).Value = 4
$ws.Cells.Item(86, 7).Value = 4
$ws.Cells.Item(86, 8).Value = 4
$ws.Cells.Item(86, 9).Value = 4
$ws.Cells.Item(87, 1).Value = "M4"
$ws.Cells.Item(87, 7).Value = 3
$ws.Cells.Item(87, 8).Value = 3
$ws.Cells.Item(87, 9).Value = 3
$ws.Cells.Item(88, 1).Value = "PF"
$ws.Cells.Item(88, 2).Value = "X"
$ws.Cells.Item(88, 3).Value = "X"
$ws.Cells.Item(88, 4).Value = "X"
$ws.Cells.Item(88, 7).Value = "X"
$ws.Cells.Item(88, 9).Value = "X"

# ---- Block: 2 chance (header row 90) ----
$ws.Range("A59:N64").Copy()
$ws.Range("A91").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(90, 1).Value = "2 chance"
$ws.Cells.Item(91, 1).Value = "Marco/Pagina"
$ws.Cells.Item(91, 2).Value = 1
$ws.Cells.Item(91, 3).Value = 2
$ws.Cells.Item(91, 4).Value = 4
$ws.Cells.Item(91, 5).Value = 2
$ws.Cells.Item(91, 6).Value = 1
$ws.Cells.Item(91, 7).Value = 3
$ws.Cells.Item(91, 8).Value = 4
$ws.Cells.Item(91, 9).Value = 5
$ws.Cells.Item(91, 10).Value = 1
$ws.Cells.Item(91, 11).Value = 6
$ws.Cells.Item(91, 12).Value = 1
$ws.Cells.Item(91, 13).Value = 2
$ws.Cells.Item(91, 14).Value = 3
$ws.Cells.Item(92, 1).Value = "M1"
$ws.Cells.Item(92, 2).Value = 1
$ws.Cells.Item(92, 3).Value = 1
$ws.Cells.Item(92, 4).Value = 1
$ws.Cells.Item(92, 5).Value = 1
$ws.Cells.Item(92, 6).Value = "1*"
$ws.Cells.Item(92, 7).Value = "1*"
$ws.Cells.Item(92, 8).Value = "1*"
$ws.Cells.Item(92, 9).Value = 1
$ws.Cells.Item(93, 1).Value = "M2"
$ws.Cells.Item(93, 3).Value = 2
$ws.Cells.Item(93, 4).Value = 2
$ws.Cells.Item(93, 5).Value = "2*"
$ws.Cells.Item(93, 6).Value = "2*"
$ws.Cells.Item(93, 7).Value = "2*"
$ws.Cells.Item(93, 8).Value = "2*"
$ws.Cells.Item(93, 9).Value = 2
$ws.Cells.Item(94, 1).Value = "M3"
$ws.Cells.Item(94, 4).Value = 4
$ws.Cells.Item(94, 5).Value = 4
$ws.Cells.Item(94, 6).Value = 4
$ws.Cells.Item(94, 7).Value = 4
$ws.Cells.Item(94, 8).Value = "4*"
$ws.Cells.Item(94, 9).Value = 4
$ws.Cells.Item(95, 1).Value = "M4"
$ws.Cells.Item(95, 7).Value = 3
$ws.Cells.Item(95, 8).Value = 3
$ws.Cells.Item(95, 9).Value = 5
$ws.Cells.Item(96, 1).Value = "PF"
$ws.Cells.Item(96, 2).Value = "X"
$ws.Cells.Item(96, 3).Value = "X"
$ws.Cells.Item(96, 4).Value = "X"
$ws.Cells.Item(96, 7).Value = "X"
$ws.Cells.Item(96, 9).Value = "X"

# ---- Block: OPT (header row 98) ----
$ws.Range("A59:N64").Copy()
$ws.Range("A99").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(98, 1).Value = "OPT"
$ws.Cells.Item(99, 1).Value = "Marco/Pagina"
$ws.Cells.Item(99, 2).Value = 1
$ws.Cells.Item(99, 3).Value = 2
$ws.Cells.Item(99, 4).Value = 4
$ws.Cells.Item(99, 5).Value = 2
$ws.Cells.Item(99, 6).Value = 1
$ws.Cells.Item(99, 7).Value = 3
$ws.Cells.Item(99, 8).Value = 4
$ws.Cells.Item(99, 9).Value = 5
$ws.Cells.Item(99, 10).Value = 1
$ws.Cells.Item(99, 11).Value = 6
$ws.Cells.Item(99, 12).Value = 1
$ws.Cells.Item(99, 13).Value = 2
$ws.Cells.Item(99, 14).Value = 3
$ws.Cells.Item(100, 1).Value = "M1"
$ws.Cells.Item(100, 2).Value = 1
$ws.Cells.Item(100, 3).Value = 1
$ws.Cells.Item(100, 4).Value = 1
$ws.Cells.Item(100, 5).Value = 1
$ws.Cells.Item(100, 6).Value = 1
$ws.Cells.Item(100, 7).Value = 1
$ws.Cells.Item(100, 8).Value = 1
$ws.Cells.Item(100, 9).Value = 1
$ws.Cells.Item(101, 1).Value = "M2"
$ws.Cells.Item(101, 3).Value = 2
$ws.Cells.Item(101, 4).Value = 2
$ws.Cells.Item(101, 5).Value = 2
$ws.Cells.Item(101, 6).Value = 2
$ws.Cells.Item(101, 7).Value = 2
$ws.Cells.Item(101, 8).Value = 2
$ws.Cells.Item(101, 9).Value = 2
$ws.Cells.Item(102, 1).Value = "M3"
$ws.Cells.Item(102, 4).Value = 4
$ws.Cells.Item(102, 5).Value = 4
$ws.Cells.Item(102, 6).Value = 4
$ws.Cells.Item(102, 7).Value = 4
$ws.Cells.Item(102, 8).Value = 4
$ws.Cells.Item(102, 9).Value = 5
$ws.Cells.Item(103, 1).Value = "M4"
$ws.Cells.Item(103, 7).Value = 3
$ws.Cells.Item(103, 8).Value = 3
$ws.Cells.Item(103, 9).Value = 3
$ws.Cells.Item(104, 1).Value = "PF"
$ws.Cells.Item(104, 2).Value = "X"
$ws.Cells.Item(104, 3).Value = "X"
$ws.Cells.Item(104, 4).Value = "X"
$ws.Cells.Item(104, 7).Value = "X"
$ws.Cells.Item(104, 9).Value = "X"

# --- Restore the view: active cell P107 (matches target selection) ---
$ws.Range("P107").Select()
